$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new addr_state dummy rows (NY, CA) and shift the rest down ---
# Current column D (rows 4-12) holds the addr_state dummies after the
# reference category in D2/D3. We need to make room for "addr_state:NY"
# (new row 4) and "addr_state:CA" (new row 6), pushing the remaining
# dummies down, and append "addr_state:IL_TX" before the last three
# existing categories.

# Read the existing addr_state dummy values (D4:D12) before overwriting them.
$oldValues = @()
for ($r = 4; $r -le 12; $r++) {
    $oldValues += $ws.Cells.Item($r, 4).Value2
}
# $oldValues now contains, in order:
# OK_TN_MO_LA_MD_NC, UT_KY_AZ_NJ, AR_MI_PA_OH_MN, RI_MA_DE_SD_IN,
# GA_WA_OR, WI_MT, IL_CT, KS_SC_CO_VT_AK_MS, WV_NH_WY_DC_ME_ID

$newOrder = @(
    "addr_state:NY",
    $oldValues[0],
    "addr_state:CA",
    $oldValues[1],
    $oldValues[2],
    $oldValues[3],
    $oldValues[4],
    $oldValues[5],
    "addr_state:IL_TX",
    $oldValues[6],
    $oldValues[7],
    $oldValues[8]
)

for ($i = 0; $i -lt $newOrder.Length; $i++) {
    $r = 4 + $i
    $ws.Cells.Item($r, 4).Value = $newOrder[$i]
}

# The previous loop only touches rows 4-15, but rows 13-15 did not exist
# before, so they do not yet carry the same "addr_state" list formatting
# (left aligned, bottom border) as the rest of the D column. Copy that
# formatting down from the template cell used throughout the list.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D13:D15").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# --- New verification_status dummy variable block in column E ---
$ws.Range("E1").Value = "verification_status"
$ws.Range("E2").Value = "verification_status:Verified"
$ws.Range("E3").Value = "verification_status:Source Verified"
$ws.Range("E4").Value = "verification_status:Not Verified"

# --- Column E width to fit the new content ---
$ws.Columns.Item(5).ColumnWidth = 28.45

# --- Restore the active selection to match the edited cell ---
$ws.Range("E9").Select() | Out-Null
